$wb = $excel.ActiveWorkbook
$validLogin = $wb.Worksheets.Item("ValidLogin")

$newSheet = $wb.Worksheets.Add($validLogin)
$newSheet.Name = "InvalidLogin"

$newSheet.Range("A1").Value = "Username"
$newSheet.Range("B1").Value = "Password"
$newSheet.Range("A2").Value = "abcd"
$newSheet.Range("B2").Value = "xyz"

$validLogin.Activate()
$validLogin.Range("A3").Select() | Out-Null

$newSheet.Activate()
$newSheet.Range("B3").Select() | Out-Null
